$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the USD amount value in T2 (568709 -> 570690)
$ws.Range("T2").Value = 570690

# Move the active selection from T2 to T3
$ws.Range("T3").Select()
